$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 9 "media" rows (IDs 10200-10208), which are rows 12-20 (1-indexed)
$ws.Range("A12:T20").EntireRow.Delete()

# Restore the selection to mirror the post-edit cursor position
$ws.Range("D5").Select()
